$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$ws = $wb.Worksheets.Item("start_price")
$ws.Range("A2").Value = 2.055

# --- Linear sheet ---
$ws = $wb.Worksheets.Item("Linear")
$ws.Range("B2").Value = -0.0002692854079815117
$ws.Range("B3").Value = -0.008157298951845382
$ws.Range("B4").Value = 0.02787177565366872
$ws.Range("B5").Value = "[1.0, 0.22552798732738358, 0.056398891485732607, 0.08560144702398881, 0.046777563525830075, 0.05498943797636096, 0.200594025492806, 0.36279221525143707, 0.21169202380457786, 0.04356483841501016, 0.0377714186123273, 0.04325980061669369, 0.06308859578436626, 0.20272374636548202, 0.3660825975623179, 0.20623095146838655, 0.019408755245307822, 0.04365438549181666, 0.037004998520316995, 0.03738080180442202]"

# --- NonLinear sheet ---
$ws = $wb.Worksheets.Item("NonLinear")
$ws.Range("B3").Value = 1.068016194331984
$ws.Range("B4").Value = 0.01055150928231232
$ws.Range("B5").Value = 0.1823669624468967
$ws.Range("B6").Value = 0.02407186804626683
$ws.Range("B7").Value = 0.003331126887409466
$ws.Range("B8").Value = -0.1216428470706532
$ws.Range("B9").Value = 0.03179532681867699
$ws.Range("B10").Value = "[1.0, 0.2259759267733879, 0.06677700088527702, 0.10164349285820642, 0.06122543810823906, 0.06935388684336177, 0.2032340672974424, 0.35332431635126543, 0.20987094283380378, 0.055226905722150325, 0.05120563627544868, 0.059493819152457744, 0.07571793970704686, 0.2020011977673216, 0.35928795187770274, 0.20778471558690997, 0.031630825573668545, 0.05512319100171853, 0.04906806307160904, 0.04777514371612177]"
